$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Un-minimize the workbook window
$wb.Windows.Item(1).WindowState = -4143  # xlNormal

# Column C: set width to fit (bestFit) - matches <col min="3" max="3" width="10" bestFit="1" customWidth="1"/>
$ws.Columns.Item(3).ColumnWidth = 9.1

# Row 2
$ws.Range("A2").Value = 2048
$ws.Range("C2").Value = 0.0009
$ws.Range("D2").Value = 81.9
$ws.Range("E2").Value = 89
$ws.Range("F2").Value = 101.1

# Row 3
$ws.Range("A3").Value = 2048
$ws.Range("C3").Value = 0.0028
$ws.Range("D3").Value = 128.3
$ws.Range("E3").Value = 137.8
$ws.Range("F3").Value = 165.6

# Row 4
$ws.Range("A4").Value = 2048
$ws.Range("C4").Value = 0.0127
$ws.Range("D4").Value = 256.5
$ws.Range("E4").Value = 269.2
$ws.Range("F4").Value = 294.4

# Row 5
$ws.Range("A5").Value = 1024
$ws.Range("C5").Value = 0.0000004
$ws.Range("D5").Value = 81.9
$ws.Range("E5").Value = 88.1
$ws.Range("F5").Value = 91

# Row 6
$ws.Range("A6").Value = 1024
$ws.Range("C6").Value = 0.0000033
$ws.Range("D6").Value = 128.7
$ws.Range("E6").Value = 138.8
$ws.Range("F6").Value = 141.2

# Row 7
$ws.Range("A7").Value = 1024
$ws.Range("C7").Value = 0.000052
$ws.Range("D7").Value = 256.8
$ws.Range("E7").Value = 275.5
$ws.Range("F7").Value = 272.8

# Rows 8, 9, 10: clear contents (keep formatting/style)
$ws.Range("A8:F10").ClearContents()

# Select C5 on the sheet (matches <selection activeCell="C5" sqref="C5"/>)
$ws.Range("C5").Select()
